$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.243.44'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.063.62'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.67%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '391.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.62%  '
$ws.Range("E7").Value = '  -2.28%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.584'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.80'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.47%  '
$ws.Range("E11").Value = '  +0.38%  '
$ws.Range("E12").Value = '  -1.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.545.52'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.28'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.65'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.57%  '
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.02'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.96%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.064.40'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.60%  '
$ws.Range("E18").Value = '  -1.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '51.253.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.86%  '
$ws.Range("E20").Value = '  +1.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0954'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.70'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '264.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.15'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.94%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.59%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.10'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.70%  '
$ws.Range("E30").Value = '  -7.74%  '
$ws.Range("E31").Value = '  -2.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.71'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0486'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '35.78'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.74%  '
$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '49.94'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.21%  '
$ws.Range("B36").Value = 'Toncoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.03'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.31'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.292'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '128.37'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.60'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.09%  '
$ws.Range("E42").Value = '  -1.99%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.81'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.20%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.115'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.50'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.67'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.47'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("E48").Value = '  -2.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.066.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.46%  '
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.883'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.41%  '
